# Refresh the cryptos price/volume snapshot (and two swapped rows: Bittensor <-> InjectiveProtocol).
# D-column cells that are plain numeric text (e.g. "5.56") are written with a leading
# apostrophe so Excel's COM layer keeps them as text (matching the workbook's original
# inline-string cells) instead of silently coercing them to numeric values; the style is
# then reset to Normal so the quote-prefix formatting does not linger on the cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "65.823.57"
$ws.Range("E2").Value = "  -2.37%  "
# Row 3
$ws.Range("D3").Value = "3.282.38"
$ws.Range("E3").Value = "  -0.86%  "
# Row 4
$ws.Range("E4").Value = "  +0.01%  "
# Row 5
$ws.Range("D5").Value = "'572.70"
$ws.Range("E5").Value = "  -0.82%  "
# Row 6
$ws.Range("D6").Value = "'177.33"
# Row 7
$ws.Range("D7").Value = "'0.631"
$ws.Range("E7").Value = "  +4.70%  "
# Row 8
$ws.Range("E8").Value = "  +0.00%  "
# Row 9
$ws.Range("E9").Value = "  -2.74%  "
# Row 10
$ws.Range("E10").Value = "  +0.86%  "
# Row 11
$ws.Range("D11").Value = "'0.399"
$ws.Range("E11").Value = "  -2.84%  "
# Row 12
$ws.Range("D12").Value = "3.853.25"
$ws.Range("E12").Value = "  -0.87%  "
# Row 13
$ws.Range("E13").Value = "  -3.80%  "
# Row 14
$ws.Range("D14").Value = "'26.53"
$ws.Range("E14").Value = "  -3.41%  "
# Row 15
$ws.Range("D15").Value = "65.889.79"
$ws.Range("E15").Value = "  -2.58%  "
# Row 16
$ws.Range("E16").Value = "  -2.63%  "
# Row 17
$ws.Range("D17").Value = "3.283.05"
$ws.Range("E17").Value = "  -0.77%  "
# Row 18
$ws.Range("D18").Value = "'436.57"
$ws.Range("E18").Value = "  -1.55%  "
# Row 19
$ws.Range("D19").Value = "'5.56"
$ws.Range("E19").Value = "  -2.44%  "
# Row 20
$ws.Range("D20").Value = "'13.17"
$ws.Range("E20").Value = "  -3.03%  "
# Row 21
$ws.Range("D21").Value = "'7.39"
$ws.Range("E21").Value = "  -4.75%  "
# Row 22
$ws.Range("D22").Value = "'72.58"
$ws.Range("E22").Value = "  -2.00%  "
# Row 23
$ws.Range("E23").Value = "  +0.13%  "
# Row 24
$ws.Range("D24").Value = "3.430.95"
$ws.Range("E24").Value = "  -0.75%  "
# Row 25
$ws.Range("D25").Value = "'0.508"
$ws.Range("E25").Value = "  -1.78%  "
# Row 26
$ws.Range("E26").Value = "  -5.12%  "
# Row 27
$ws.Range("D27").Value = "'0.195"
$ws.Range("E27").Value = "  +3.69%  "
# Row 28
$ws.Range("D28").Value = "'8.86"
$ws.Range("E28").Value = "  -2.05%  "
# Row 29
$ws.Range("E29").Value = "  +0.11%  "
# Row 30
$ws.Range("E30").Value = "  -2.24%  "
# Row 31
$ws.Range("D31").Value = "'22.28"
$ws.Range("E31").Value = "  -2.92%  "
# Row 32
$ws.Range("E32").Value = "  +0.08%  "
# Row 33
$ws.Range("D33").Value = "'5.14"
$ws.Range("E33").Value = "  -3.75%  "
# Row 34
$ws.Range("D34").Value = "'6.59"
$ws.Range("E34").Value = "  -3.27%  "
# Row 35
$ws.Range("D35").Value = "'1.18"
$ws.Range("E35").Value = "  -4.77%  "
# Row 36
$ws.Range("D36").Value = "'158.72"
$ws.Range("E36").Value = "  -2.50%  "
# Row 37
$ws.Range("D37").Value = "'1.44"
$ws.Range("E37").Value = "  -4.90%  "
# Row 38
$ws.Range("D38").Value = "'26.66"
$ws.Range("E38").Value = "  -1.93%  "
# Row 39
$ws.Range("D39").Value = "'1.78"
$ws.Range("E39").Value = "  -3.92%  "
# Row 40
$ws.Range("D40").Value = "2.772.79"
$ws.Range("E40").Value = "  +0.59%  "
# Row 41
$ws.Range("D41").Value = "'0.778"
$ws.Range("E41").Value = "  -1.62%  "
# Row 42
$ws.Range("D42").Value = "'4.31"
$ws.Range("E42").Value = "  -3.91%  "
# Row 43
$ws.Range("D43").Value = "'40.30"
$ws.Range("E43").Value = "  +0.28%  "
# Row 44
$ws.Range("D44").Value = "'6.03"
$ws.Range("E44").Value = "  -3.42%  "
# Row 45
$ws.Range("D45").Value = "'0.0655"
$ws.Range("E45").Value = "  -2.51%  "
# Row 46
$ws.Range("D46").Value = "'2.28"
$ws.Range("E46").Value = "  -5.29%  "
# Row 47
$ws.Range("B47").Value = "InjectiveProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D47").Value = "'23.33"
$ws.Range("E47").Value = "  -6.04%  "
# Row 48
$ws.Range("B48").Value = "Bittensor"
$ws.Range("C48").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D48").Value = "'316.87"
$ws.Range("E48").Value = "  -3.24%  "
# Row 49
$ws.Range("D49").Value = "'0.0268"
$ws.Range("E49").Value = "  -1.92%  "
# Row 50
$ws.Range("E50").Value = "  +2.43%  "
# Row 51
$ws.Range("E51").Value = "  +0.05%  "

# Strip the quote-prefix style picked up above so affected D-column cells keep the
# workbook default (unstyled) appearance.
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Style = "Normal"
$ws.Range("D7").Style = "Normal"
$ws.Range("D11").Style = "Normal"
$ws.Range("D14").Style = "Normal"
$ws.Range("D18").Style = "Normal"
$ws.Range("D19").Style = "Normal"
$ws.Range("D20").Style = "Normal"
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").Style = "Normal"
$ws.Range("D25").Style = "Normal"
$ws.Range("D27").Style = "Normal"
$ws.Range("D28").Style = "Normal"
$ws.Range("D31").Style = "Normal"
$ws.Range("D33").Style = "Normal"
$ws.Range("D34").Style = "Normal"
$ws.Range("D35").Style = "Normal"
$ws.Range("D36").Style = "Normal"
$ws.Range("D37").Style = "Normal"
$ws.Range("D38").Style = "Normal"
$ws.Range("D39").Style = "Normal"
$ws.Range("D41").Style = "Normal"
$ws.Range("D42").Style = "Normal"
$ws.Range("D43").Style = "Normal"
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").Style = "Normal"
$ws.Range("D46").Style = "Normal"
$ws.Range("D47").Style = "Normal"
$ws.Range("D48").Style = "Normal"
$ws.Range("D49").Style = "Normal"
